$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Item(30, 1).Value = "hybrid_cbf_cf_w0.13cf_w0.87cbf_popularity1000_biasGiusti"
$ws.Cells.Item(30, 2).NumberFormat = "@"
$ws.Cells.Item(30, 2).Value = "0.00556"
